# feat: add 2022-Q1 data
#
# 1. Insert a new "2022-Q1" sheet (same 8-column layout as "2021-Q3") right
#    before the "总计" sheet.
# 2. Insert a new top row in "总计" summarizing the 2022-Q1 quarter, pushing
#    the existing rows down.

$wb = $excel.ActiveWorkbook

$wsTotalOld = $wb.Worksheets.Item(3)   # "总计" (before insertion)
$wsRef      = $wb.Worksheets.Item(2)   # "2021-Q3" - used as a template since
                                        # it already has the 8-column layout.

# Duplicate the template sheet and drop it in right before "总计" so the tab
# order ends up 2021-Q2, 2021-Q3, 2022-Q1, 总计.
$wsRef.Copy($wsTotalOld)
$newSheet = $wb.Worksheets.Item(3)
$newSheet.Name = "2022-Q1"

# "总计" shifted from index 3 to index 4 - fetch it again by name.
$wsTotal = $wb.Worksheets.Item("总计")

# Columns B and D:G hold numeric-looking text ("011164", "84.49", ...); force
# text format first so Excel doesn't coerce them into numbers (which would
# drop the leading zeros / change the stored type).
$newSheet.Range("B2:B5").NumberFormat = "@"
$newSheet.Range("D2:G5").NumberFormat = "@"

$newSheet.Range("B1").Value = "基金代码"
$newSheet.Range("C1").Value = "基金名称"
$newSheet.Range("D1").Value = "基金规模"
$newSheet.Range("E1").Value = "股票总仓位"
$newSheet.Range("F1").Value = "仓位占比"
$newSheet.Range("G1").Value = "持有市值(亿元)"
$newSheet.Range("H1").Value = "仓位排名"

$newSheet.Range("A2").Value = 0
$newSheet.Range("B2").Value = "011164"
$newSheet.Range("C2").Value = "富国兴远优选12个月持有期混合A"
$newSheet.Range("D2").Value = "84.49"
$newSheet.Range("E2").Value = "82.18"
$newSheet.Range("F2").Value = "3.25"
$newSheet.Range("G2").Value = "2.7459"
$newSheet.Range("H2").Value = 8

$newSheet.Range("A3").Value = 1
$newSheet.Range("B3").Value = "001186"
$newSheet.Range("C3").Value = "富国文体健康股票A"
$newSheet.Range("D3").Value = "18.44"
$newSheet.Range("E3").Value = "86.61"
$newSheet.Range("F3").Value = "3.84"
$newSheet.Range("G3").Value = "0.7081"
$newSheet.Range("H3").Value = 7

$newSheet.Range("A4").Value = 2
$newSheet.Range("B4").Value = "011165"
$newSheet.Range("C4").Value = "富国兴远优选12个月持有期混合C"
$newSheet.Range("D4").Value = "18.89"
$newSheet.Range("E4").Value = "82.18"
$newSheet.Range("F4").Value = "3.25"
$newSheet.Range("G4").Value = "0.6139"
$newSheet.Range("H4").Value = 8

$newSheet.Range("A5").Value = 3
$newSheet.Range("B5").Value = "011125"
$newSheet.Range("C5").Value = "富国文体健康股票C"
$newSheet.Range("D5").Value = "0.71"
$newSheet.Range("E5").Value = "86.61"
$newSheet.Range("F5").Value = "3.84"
$newSheet.Range("G5").Value = "0.0273"
$newSheet.Range("H5").Value = 7

# --- Update "总计": shift existing rows down one and write the new
#     2022-Q1 summary row at the top (row 2). ---
$wsTotal.Range("A4").Value = 2
$wsTotal.Range("B4").Value = "2021-Q2"
$wsTotal.Range("C4").Value = 2
$wsTotal.Range("D4").Value = 3.02

$wsTotal.Range("A3").Value = 1
$wsTotal.Range("B3").Value = "2021-Q3"
$wsTotal.Range("C3").Value = 4
$wsTotal.Range("D3").Value = 3.71

# A2 already carries the bold/bordered "index column" style - copy it onto
# the new A4 before overwriting A2's own value.
$wsTotal.Range("A2").Copy($wsTotal.Range("A4"))
$wsTotal.Range("A4").Value = 2

$wsTotal.Range("A2").Value = 0
$wsTotal.Range("B2").Value = "2022-Q1"
$wsTotal.Range("C2").Value = 4
$wsTotal.Range("D2").Value = 4.1

# Restore the original active sheet/tab selection (2021-Q2).
$wb.Worksheets.Item(1).Activate()
